$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.114.49"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "2.253.04"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'307.30"
$ws.Range("E5").Value = "  -4.87%  "

$ws.Range("D6").Value = "'98.74"
$ws.Range("E6").Value = "  -3.00%  "

$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  -3.79%  "

$ws.Range("D10").Value = "'35.67"
$ws.Range("E10").Value = "  -4.75%  "

$ws.Range("D12").Value = "'7.32"
$ws.Range("E12").Value = "  -5.06%  "

$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("D14").Value = "2.595.73"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.257.82"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.839"
$ws.Range("E16").Value = "  -2.30%  "

$ws.Range("D17").Value = "'13.87"
$ws.Range("E17").Value = "  -2.22%  "

$ws.Range("D18").Value = "44.009.77"
$ws.Range("E18").Value = "  +0.72%  "

$ws.Range("D19").Value = "'12.79"
$ws.Range("E19").Value = "  -7.54%  "

$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").Value = "'6.35"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").Value = "'65.44"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'241.45"
$ws.Range("E23").Value = "  +2.10%  "

$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  -7.08%  "

$ws.Range("E25").Value = "  -8.55%  "

$ws.Range("D26").Value = "'0.990"
$ws.Range("E26").Value = "  -1.22%  "

$ws.Range("D27").Value = "'10.14"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").Value = "'37.51"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  -2.61%  "

$ws.Range("D30").Value = "'6.18"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").Value = "'20.10"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("D32").Value = "'157.04"
$ws.Range("E32").Value = "  -1.98%  "

$ws.Range("D33").Value = "'3.52"
$ws.Range("E33").Value = "  +10.00%  "

$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("E37").Value = "  -5.91%  "

$ws.Range("D38").Value = "'1.86"
$ws.Range("E38").Value = "  -4.15%  "

$ws.Range("D39").Value = "'15.54"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").Value = "'3.89"
$ws.Range("E40").Value = "  -9.62%  "

$ws.Range("D41").Value = "'3.39"
$ws.Range("E41").Value = "  -11.17%  "

$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("D44").Value = "1.766.22"
$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("D45").Value = "'87.50"
$ws.Range("E45").Value = "  +5.54%  "

$ws.Range("E46").Value = "  -3.93%  "

$ws.Range("D47").Value = "'5.13"
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").Value = "'101.61"
$ws.Range("E48").Value = "  -2.25%  "

$ws.Range("D49").Value = "'8.25"
$ws.Range("E49").Value = "  -2.58%  "

$ws.Range("D50").Value = "'70.31"
$ws.Range("E50").Value = "  -5.81%  "

$ws.Range("D51").Value = "'55.57"
$ws.Range("E51").Value = "  -5.82%  "
